# "data till 7Feb 1PM" -- enter newly-collected payments (columns K = 5 Feb,
# L = 6 Feb) for a batch of retailers on the daily-collection sheet.
# Row "Total" (F) and the sheet-wide date totals (row 2) are formula driven
# and recalculate automatically once the underlying figures are entered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain entries (no special highlight fill) -------------------------------
$plainEntries = @(
    @{Cell = "K5";  Value = 2000},
    @{Cell = "K10"; Value = 5000},
    @{Cell = "K12"; Value = 1000},
    @{Cell = "K14"; Value = 2000},
    @{Cell = "K19"; Value = 2000},
    @{Cell = "L20"; Value = 1000},
    @{Cell = "K23"; Value = 3000},
    @{Cell = "K31"; Value = 1000},
    @{Cell = "L31"; Value = 500},
    @{Cell = "K40"; Value = 2000},
    @{Cell = "K41"; Value = 3000},
    @{Cell = "L45"; Value = 3000},
    @{Cell = "K46"; Value = 900},
    @{Cell = "K48"; Value = 2000},
    @{Cell = "L48"; Value = 3000},
    @{Cell = "K50"; Value = 2000},
    @{Cell = "L52"; Value = 1000},
    @{Cell = "L53"; Value = 2000},
    @{Cell = "L72"; Value = 1500}
)

foreach ($entry in $plainEntries) {
    $ws.Range($entry.Cell).Value = $entry.Value
}

# Entries that also pick up the purple "highlighted" fill (same shade the
# author already used elsewhere, e.g. G3, H35, J18, ...) -------------------
$highlightColor = 0x93697B   # OLE BGR for RGB 7B6993

$highlightedEntries = @(
    @{Cell = "L18"; Value = 5000},
    @{Cell = "K26"; Value = 1000},
    @{Cell = "L26"; Value = 1000},
    @{Cell = "L30"; Value = 10000},
    @{Cell = "L35"; Value = 3000},
    @{Cell = "L47"; Value = 2000},
    @{Cell = "L71"; Value = 6000},
    @{Cell = "L81"; Value = 1000}
)

foreach ($entry in $highlightedEntries) {
    $rng = $ws.Range($entry.Cell)
    $rng.Value = $entry.Value
    $rng.Interior.Color = $highlightColor
}

# New retailer area label that was missing before -------------------------
$ws.Range("D86").Value = "KALER"

# Move the active selection to L3, where the user was last working --------
$ws.Range("L3").Select()
